$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings_recources")

# Defaultní přípona exportu (Katalog): xlsm -> xlsx
$ws.Range("B28").Value = "xlsx"

# nastavení zoomu celé aplikace (default: 100 %): 90 -> 80
# (force text storage, matching the sheet's existing "numbers stored as text" convention,
# then drop back to the default/Normal style so no stray cell formatting is introduced)
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "80"
$ws.Range("B30").Style = "Normal"

# nový řádek 32: precise
$ws.Range("B32").Value = "precise"
